# TC6_SearchResults_Typeahead.xlsx - "Changes for New UI Prod"
#
# Sheet1 (TC6_SearchResults_Typeahead): insert two new keyword rows
#   (CLICK_PRE_ENTERTEXT before ENTERTEXT, and WAIT after both
#   ENTERTEXT and PRESS_ENTER) and update the selected cell.
# Sheet2 (Testdata): append two new data rows (EleType1/EleType2 -> JSElement).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

function Set-ThinBorder($rng) {
    $rng.Borders.LineStyle = 1
    $rng.Borders.Weight = 2
}

# ---------------------------------------------------------------
# Sheet1: TC6_SearchResults_Typeahead
# ---------------------------------------------------------------
# Current layout (rows 1-6):
#   1 Header
#   2 TC6_SearchResults_Typeahead | NAVIGATE_URL | | | Baseurl
#   3 | ENTERTEXT | SearchBoxHomePage | CSS | TypeaheadText
#   4 | VERIFY_WEBELEMENT_PRESENT | Typeahead | CSS | Typeahead
#   5 | PRESS_ENTER | SearchBoxHomePage | CSS |
#   6 | VERIFY_TEXT_PRESENT | ValidSearchHeader1 | CSS | validSearchText
#
# Target layout (rows 1-9):
#   1 Header
#   2 TC6_SearchResults_Typeahead | NAVIGATE_URL | | | Baseurl
#   3 | CLICK_PRE_ENTERTEXT | SearchBoxHomePage | CSS |        (NEW)
#   4 | ENTERTEXT | SearchBoxHomePage | CSS | TypeaheadText
#   5 | WAIT | | |                                              (NEW)
#   6 | VERIFY_WEBELEMENT_PRESENT | Typeahead | CSS | Typeahead
#   7 | PRESS_ENTER | SearchBoxHomePage | CSS |
#   8 | WAIT | | |                                               (NEW)
#   9 | VERIFY_TEXT_PRESENT | ValidSearchHeader1 | CSS | validSearchText

# Insert a new row above the current row 3 (ENTERTEXT) for CLICK_PRE_ENTERTEXT.
$ws1.Range("A3:E3").EntireRow.Insert()

# Insert a new row above the current row 5 (which, after the previous insert,
# is PRESS_ENTER) for the WAIT row that sits between VERIFY_WEBELEMENT_PRESENT
# and PRESS_ENTER -- i.e. directly after ENTERTEXT.
$ws1.Range("A5:E5").EntireRow.Insert()

# Insert a new row above the current row 8 (VERIFY_TEXT_PRESENT) for the
# second WAIT row, which sits right after PRESS_ENTER.
$ws1.Range("A8:E8").EntireRow.Insert()

# --- Fill in the new row 3: CLICK_PRE_ENTERTEXT ---
$ws1.Range("B3").Value = "CLICK_PRE_ENTERTEXT"
$ws1.Range("C3").Value = "SearchBoxHomePage"
$ws1.Range("D3").Value = "CSS"
Set-ThinBorder($ws1.Range("A3:E3"))

# --- Fill in the new row 5: WAIT ---
$ws1.Range("B5").Value = "WAIT"
Set-ThinBorder($ws1.Range("A5:E5"))

# --- Fill in the new row 8: WAIT ---
$ws1.Range("B8").Value = "WAIT"
Set-ThinBorder($ws1.Range("A8:E8"))

# Update the selected cell shown when the sheet is opened.
$ws1.Range("B5").Select()

# ---------------------------------------------------------------
# Sheet2: Testdata
# ---------------------------------------------------------------
# Append two new rows at the bottom of the data table:
#   8 EleType1 | JSElement
#   9 EleType2 | JSElement
$ws2.Range("A8").Value = "EleType1"
$ws2.Range("B8").Value = "JSElement"
Set-ThinBorder($ws2.Range("A8:B8"))

$ws2.Range("A9").Value = "EleType2"
$ws2.Range("B9").Value = "JSElement"
Set-ThinBorder($ws2.Range("A9:B9"))

Write-Host "Edit complete"
